$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '66.955.06'
$ws.Range('E2').Value = '  +0.01%  '

# Row 3
$ws.Range('D3').Value = '3.091.56'
$ws.Range('E3').Value = '  +0.25%  '

# Row 4
$ws.Range('E4').Value = '  -0.06%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '572.08'
$ws.Range('E5').Value = '  -1.20%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '177.58'
$ws.Range('E6').Value = '  +4.63%  '

# Row 7
$ws.Range('E7').Value = '  -0.10%  '

# Row 8
$ws.Range('D8').Value = '3.088.97'
$ws.Range('E8').Value = '  +0.27%  '

# Row 9
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.512'
$ws.Range('E9').Value = '  -0.51%  '

# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '6.38'
$ws.Range('E10').Value = '  -0.46%  '

# Row 11
$ws.Range('E11').Value = '  +0.19%  '

# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.467'
$ws.Range('E12').Value = '  -1.33%  '

# Row 13
$ws.Range('E13').Value = '  -1.11%  '

# Row 14
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '35.99'
$ws.Range('E14').Value = '  -0.72%  '

# Row 15
$ws.Range('E15').Value = '  +0.66%  '

# Row 16
$ws.Range('D16').Value = '3.607.81'
$ws.Range('E16').Value = '  +0.23%  '

# Row 17
$ws.Range('D17').Value = '66.971.89'
$ws.Range('E17').Value = '  +0.05%  '

# Row 18
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '7.01'
$ws.Range('E18').Value = '  -0.34%  '

# Row 19
$ws.Range('B19').Value = 'WrappedEther'
$ws.Range('C19').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D19').Value = '3.093.76'
$ws.Range('E19').Value = '  +0.42%  '

# Row 20
$ws.Range('B20').Value = 'Chainlink'
$ws.Range('C20').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '16.56'
$ws.Range('E20').Value = '  +0.81%  '

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '484.36'
$ws.Range('E21').Value = '  +0.25%  '

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '7.68'
$ws.Range('E22').Value = '  -0.52%  '

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.684'
$ws.Range('E23').Value = '  -0.98%  '

# Row 24
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '83.35'
$ws.Range('E24').Value = '  +0.20%  '

# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.26'
$ws.Range('E25').Value = '  +0.41%  '

# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '12.57'
$ws.Range('E26').Value = '  -2.56%  '

# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '10.17'
$ws.Range('E27').Value = '  -2.33%  '

# Row 28
$ws.Range('E28').Value = '  -0.01%  '

# Row 29
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '7.86'
$ws.Range('E29').Value = '  +2.14%  '

# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '2.29'
$ws.Range('E30').Value = '  -1.24%  '

# Row 31
$ws.Range('E31').Value = '  -1.96%  '

# Row 32
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '28.02'
$ws.Range('E32').Value = '  +0.21%  '

# Row 33
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.112'
$ws.Range('E33').Value = '  +0.19%  '

# Row 34
$ws.Range('D34').Value = '0.0₃0940'
$ws.Range('E34').Value = '  +1.50%  '

# Row 35
$ws.Range('E35').Value = '  -0.01%  '

# Row 36
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '47.32'
$ws.Range('E36').Value = '  +1.96%  '

# Row 37
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '5.56'
$ws.Range('E37').Value = '  -2.56%  '

# Row 38
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.943'
$ws.Range('E38').Value = '  -1.49%  '

# Row 39
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.311'
$ws.Range('E39').Value = '  +2.88%  '

# Row 40
$ws.Range('E40').Value = '  +1.52%  '

# Row 41
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '49.05'
$ws.Range('E41').Value = '  -0.98%  '

# Row 42
$ws.Range('E42').Value = '  +0.20%  '

# Row 43
$ws.Range('B43').Value = 'dogwifhat'
$ws.Range('C43').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '2.71'
$ws.Range('E43').Value = '  +5.43%  '

# Row 44
$ws.Range('B44').Value = 'Cosmos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '8.22'
$ws.Range('E44').Value = '  -1.49%  '

# Row 45
$ws.Range('D45').Value = '2.801.64'
$ws.Range('E45').Value = '  +0.87%  '

# Row 46
$ws.Range('B46').Value = 'Monero'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '135.96'
$ws.Range('E46').Value = '  +0.68%  '

# Row 47
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '368.26'
$ws.Range('E47').Value = '  -2.58%  '

# Row 48
$ws.Range('B48').Value = 'VeChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0344'
$ws.Range('E48').Value = '  -0.71%  '

# Row 50
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '25.55'
$ws.Range('E50').Value = '  +4.05%  '

# Row 51
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.28'
$ws.Range('E51').Value = '  +5.90%  '
